$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.163.78"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.669.47"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.67"
$ws.Range("E5").Value = "  -2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5237"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2624"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06338"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.22"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07539"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "1.676.41"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.449"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5507"
$ws.Range("E14").Value = "  -4.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "66.55"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000007973"
$ws.Range("E16").Value = "  -4.30%  "
$ws.Range("D17").Value = "26.158.65"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.766"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.00"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.204"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.59"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.512"
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.86"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06370"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.352"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.525"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.414"
$ws.Range("E32").Value = "  -4.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.006"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6023"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.405"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.758"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.149"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.112.69"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01615"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8638"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.38"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.53"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.052"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05233"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.930"
$ws.Range("E51").Value = "  -1.41%  "
